# #5: fund, bonds, otherbonds, antique done
#
# 1. Rebuild "基金受益憑證" (fund) sheet: add a proper header row (it used to
#    duplicate row 2's data), add a "dealer" column, and append the common
#    metadata columns (property_category, category, date, legislator_name,
#    legislator_id, source_file, index) that every other property sheet
#    already has.
# 2. Delete "其他有價證券" (antique / other securities) sheet entirely - it
#    only ever held garbled placeholder text, never real data.
# 3. "保險" and "債務" sheets are otherwise untouched; they simply shift up
#    one position once the sheet above them is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

function Set-Header($cellRef, $text) {
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
    $ws.Range($cellRef).Value = $text
}

function Set-Data($cellRef, $value) {
    $ws.Range("B2").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
    $ws.Range($cellRef).Value = $value
}

# "2011-12-28" must stay literal text, not get auto-parsed into a date
# serial number - force text format on that one column before assigning.
function Set-TextData($cellRef, $value) {
    $ws.Range("B2").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

# ---- Header row (row 1) ----
# Columns B:H already exist but previously held stray duplicated data instead
# of labels; overwrite with the real header text.
Set-Header "B1" "name"
Set-Header "C1" "owner"
Set-Header "D1" "dealer"
Set-Header "E1" "quantity"
Set-Header "F1" "face_value"
Set-Header "G1" "currency"
Set-Header "H1" "total"
Set-Header "I1" "property_category"
Set-Header "J1" "category"
Set-Header "K1" "date"
Set-Header "L1" "legislator_name"
Set-Header "M1" "legislator_id"
Set-Header "N1" "source_file"
Set-Header "O1" "index"

# ---- Data rows (rows 2-6) ----
$rows = @(
    @{ r = 2;  A = 99;  B = "摩根富林明龍揚基金";   C = "王〇聿"; D = "屏東市林毐郵局"; E = 3927.7;   F = 17.95; G = "新臺幣"; H = 70502;  O = 99 },
    @{ r = 3;  A = 100; B = "摩根富林明全球a基金";  C = "王〇聿"; D = "屏東市林森郵局"; E = 7806.4;   F = 9.68;  G = "新臺幣"; H = 75566;  O = 100 },
    @{ r = 4;  A = 101; B = "霸菱全球新興市場基金"; C = "周麗容"; D = "第一銀行";       E = 313.791;  F = 33.7;  G = "美金";   H = 342305; O = 101 },
    @{ r = 5;  A = 102; B = "貝萊德世界能源";       C = "周麗容"; D = "第一銀行";       E = 454.17;   F = 25.78; G = "美金";   H = 359451; O = 102 },
    @{ r = 6;  A = 104; B = "富達亞高收";           C = "周麗容"; D = "第一銀行";       E = 621.3;    F = 31.87; G = "美金";   H = 19801;  O = 104 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H

    Set-Data "I$r" "fund"
    Set-Data "J$r" "normal"
    Set-TextData "K$r" "2011-12-28"
    Set-Data "L$r" "王進士"
    Set-Data "M$r" 1701
    Set-Data "N$r" "tmpf41"
    Set-Data "O$r" $row.O
}

# ---- Drop the "其他有價證券" sheet entirely ----
$wb.Worksheets.Item("其他有價證券").Delete() | Out-Null
